# "Create ExcelCreator with saving data to sheet"
#
# The template originally ships with a single worksheet ("Лист1"). This
# change adds a second worksheet ("Лист2") to the workbook that is an exact
# duplicate of the first one (same layout, merged cells, styles and
# placeholder text) and leaves it as the active / selected sheet - exactly
# what Excel itself produces when you right-click a sheet tab and choose
# "Move or Copy... -> Create a copy".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Select C4:D6 on the source sheet first; the copy below carries the same
# selection over to the new sheet.
$ws1.Range("C4:D6").Select()

# Duplicate "Лист1", inserting the copy immediately after it. This preserves
# all data, number/shared-string references, cell styles and merged cells.
$ws1.Copy($null, $ws1)

# The freshly inserted copy is now the second sheet; rename it to "Лист2".
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Лист2"

# Make sure the new sheet ends up active/selected with the C4:D6 selection,
# matching the workbook's saved view state.
$ws2.Select()
$ws2.Range("C4:D6").Select()
